$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.272.52'
$ws.Range("E2").Value = '  +0.59%  '

$ws.Range("D3").Value = '1.864.93'
$ws.Range("E3").Value = '  +0.31%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9999'
$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.94'
$ws.Range("E5").Value = '  +1.64%  '

$ws.Range("E6").Value = '  +0.06%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4682'
$ws.Range("E7").Value = '  +0.57%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2862'
$ws.Range("E8").Value = '  +2.12%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06541'

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.38'
$ws.Range("E10").Value = '  +14.80%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07910'
$ws.Range("E11").Value = '  +1.21%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '97.88'
$ws.Range("E12").Value = '  +1.64%  '

$ws.Range("D13").Value = '1.868.57'
$ws.Range("E13").Value = '  +0.52%  '

$ws.Range("E14").Value = '  +1.13%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6824'
$ws.Range("E15").Value = '  +2.90%  '

$ws.Range("E16").Value = '  -0.28%  '

$ws.Range("D17").Value = '30.271.93'
$ws.Range("E17").Value = '  +0.48%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.69'
$ws.Range("E18").Value = '  +9.14%  '

$ws.Range("E19").Value = '  -0.07%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.407'
$ws.Range("E20").Value = '  -1.33%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.000007342'
$ws.Range("E21").Value = '  +1.81%  '

$ws.Range("D22").Value = '2.113.02'
$ws.Range("E22").Value = '  +0.70%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9998'
$ws.Range("E23").Value = '  +0.00%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.174'
$ws.Range("E24").Value = '  +1.02%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '167.86'
$ws.Range("E25").Value = '  +0.87%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.268'
$ws.Range("E26").Value = '  -0.13%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.15'
$ws.Range("E27").Value = '  +1.81%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.941'
$ws.Range("E28").Value = '  +2.21%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.382'
$ws.Range("E29").Value = '  +3.90%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09858'
$ws.Range("E30").Value = '  +3.53%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.399'
$ws.Range("E31").Value = '  -0.51%  '

$ws.Range("E32").Value = '  +1.33%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.070'
$ws.Range("E33").Value = '  -0.21%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04747'
$ws.Range("E34").Value = '  +2.41%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.140'
$ws.Range("E35").Value = '  +4.48%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7119'
$ws.Range("E36").Value = '  +1.83%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.708'
$ws.Range("E37").Value = '  +0.24%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01877'
$ws.Range("E38").Value = '  +1.95%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.619'
$ws.Range("E39").Value = '  +4.29%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '76.91'
$ws.Range("E40").Value = '  +6.07%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.307'
$ws.Range("E41").Value = '  +0.60%  '

$ws.Range("E42").Value = '  +2.51%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8515'
$ws.Range("E43").Value = '  -0.09%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4191'
$ws.Range("E44").Value = '  +1.32%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9998'
$ws.Range("E45").Value = '  -0.01%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '103.61'
$ws.Range("E46").Value = '  -0.15%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '964.19'
$ws.Range("E47").Value = '  -3.68%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.232'
$ws.Range("E48").Value = '  +0.86%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.360'
$ws.Range("E49").Value = '  +2.12%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '34.25'
$ws.Range("E50").Value = '  +0.89%  '

$ws.Range("E51").Value = '  +0.49%  '
